$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 360) holds a date value ("Förändrad") that was
# bumped by one day (2023-10-03 -> 2023-10-04, serial 45202 -> 45203).
$ws.Range("C2:C360").Value = 45203
